# fix: set nullable field password (import partner)
#
# Sheet1!F3 held "=F2" (cached value "password"), copying the Password
# column from the row above. The field is nullable, so this row's
# Password should be left blank instead of auto-filled.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Drop the formula/value from F3, leaving the cell blank (style is kept).
$ws.Range("F3").ClearContents()

# Reflect where the user's view/selection ended up after the edit.
$ws.Range("F3").Select()
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
